# Populate the "companyname" column (R) from the "company" column (Q)
# for rows where companyname is currently empty, as apollo source data
# was missing the companyname field needed for mautic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 17).End(-4162).Row  # xlUp = -4162, column Q = 17

for ($r = 2; $r -le $lastRow; $r++) {
    $companyCell = $ws.Cells.Item($r, 17)      # Q: company
    $companyNameCell = $ws.Cells.Item($r, 18)  # R: companyname

    $companyValue = $companyCell.Value2
    $companyNameValue = $companyNameCell.Value2

    if (($null -ne $companyValue) -and ($companyValue -ne "")) {
        if (($null -eq $companyNameValue) -or ($companyNameValue -eq "")) {
            $companyNameCell.Value = $companyValue
        }
    }
}
